# Applies the "Updated cryptos list" price/volume refresh to Sheet1.
# Numeric-looking Price (column D) strings get a leading "'" so Excel
# keeps them as text (matching the workbook's inlineStr storage) instead
# of silently coercing them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '29.108.75'
$ws.Range('E2').Value = '  -1.68%  '
# Row 3
$ws.Range('D3').Value = '1.836.07'
$ws.Range('E3').Value = '  -1.30%  '
# Row 4
$ws.Range('D4').Value = '''0.9988'
$ws.Range('E4').Value = '  -0.08%  '
# Row 5
$ws.Range('E5').Value = '  -2.12%  '
# Row 6
$ws.Range('D6').Value = '''0.6811'
$ws.Range('E6').Value = '  -2.47%  '
# Row 7
$ws.Range('D7').Value = '''0.9995'
$ws.Range('E7').Value = '  -0.07%  '
# Row 8
$ws.Range('E8').Value = '  -2.47%  '
# Row 9
$ws.Range('D9').Value = '''0.07462'
$ws.Range('E9').Value = '  -3.24%  '
# Row 10
$ws.Range('E10').Value = '  -2.23%  '
# Row 11
$ws.Range('E11').Value = '  -1.27%  '
# Row 12
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').Value = '''5.015'
$ws.Range('E12').Value = '  -2.88%  '
# Row 13
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.793.92'
$ws.Range('E13').Value = '  -3.53%  '
# Row 14
$ws.Range('D14').Value = '''0.6769'
$ws.Range('E14').Value = '  -2.17%  '
# Row 15
$ws.Range('D15').Value = '''86.57'
$ws.Range('E15').Value = '  -6.24%  '
# Row 16
$ws.Range('D16').Value = '''6.152'
$ws.Range('E16').Value = '  -6.35%  '
# Row 17
$ws.Range('D17').Value = '29.133.37'
$ws.Range('E17').Value = '  -1.55%  '
# Row 19
$ws.Range('D19').Value = '2.079.43'
$ws.Range('E19').Value = '  -1.10%  '
# Row 20
$ws.Range('D20').Value = '''228.35'
$ws.Range('E20').Value = '  -5.58%  '
# Row 21
$ws.Range('D21').Value = '''12.49'
$ws.Range('E21').Value = '  -2.14%  '
# Row 22
$ws.Range('D22').Value = '''0.9984'
$ws.Range('E22').Value = '  -0.15%  '
# Row 23
$ws.Range('D23').Value = '''7.361'
$ws.Range('E23').Value = '  -3.34%  '
# Row 24
$ws.Range('D24').Value = '''0.9996'
$ws.Range('E24').Value = '  -0.07%  '
# Row 25
$ws.Range('D25').Value = '''160.91'
$ws.Range('E25').Value = '  +0.95%  '
# Row 26
$ws.Range('D26').Value = '''0.1439'
$ws.Range('E26').Value = '  -4.26%  '
# Row 27
$ws.Range('D27').Value = '''8.709'
$ws.Range('E27').Value = '  -2.32%  '
# Row 28
$ws.Range('E28').Value = '  -1.46%  '
# Row 29
$ws.Range('D29').Value = '''1.501'
$ws.Range('E29').Value = '  -2.33%  '
# Row 30
$ws.Range('D30').Value = '''4.249'
$ws.Range('E30').Value = '  -0.03%  '
# Row 31
$ws.Range('E31').Value = '  -1.21%  '
# Row 32
$ws.Range('B32').Value = 'Toncoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D32').Value = '''1.196'
$ws.Range('E32').Value = '  +0.25%  '
# Row 33
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = '''0.05409'
$ws.Range('E33').Value = '  +6.25%  '
# Row 34
$ws.Range('D34').Value = '''0.7525'
$ws.Range('E34').Value = '  -3.08%  '
# Row 35
$ws.Range('D35').Value = '''1.855'
$ws.Range('E35').Value = '  -2.12%  '
# Row 36
$ws.Range('D36').Value = '''1.128'
$ws.Range('E36').Value = '  -2.14%  '
# Row 37
$ws.Range('E37').Value = '  -0.17%  '
# Row 38
$ws.Range('D38').Value = '1.304.58'
$ws.Range('E38').Value = '  -1.49%  '
# Row 39
$ws.Range('D39').Value = '''0.01813'
$ws.Range('E39').Value = '  -3.31%  '
# Row 40
$ws.Range('D40').Value = '''2.717'
$ws.Range('E40').Value = '  -0.60%  '
# Row 41
$ws.Range('D41').Value = '''0.9364'
$ws.Range('E41').Value = '  -2.32%  '
# Row 42
$ws.Range('D42').Value = '''6.051'
$ws.Range('E42').Value = '  +4.05%  '
# Row 43
$ws.Range('D43').Value = '''0.08491'
$ws.Range('E43').Value = '  +33.18%  '
# Row 44
$ws.Range('D44').Value = '''104.89'
$ws.Range('E44').Value = '  -1.46%  '
# Row 46
$ws.Range('D46').Value = '2.006.87'
$ws.Range('E46').Value = '  +0.11%  '
# Row 47
$ws.Range('D47').Value = '''0.5177'
$ws.Range('E47').Value = '  -0.74%  '
# Row 48
$ws.Range('E48').Value = '  -3.09%  '
# Row 49
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = '''9.427'
$ws.Range('E49').Value = '  -3.55%  '
# Row 50
$ws.Range('E50').Value = '  -0.97%  '
# Row 51
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Value = '''63.64'
$ws.Range('E51').Value = '  -1.44%  '
